# Update the "Opt Portfolio with View" column (D) values for rows 2-8
# on the active sheet, per the new computed portfolio weights.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.0004904227247809496
$ws.Range("D3").Value = [double]"1.923671547856839e-17"
$ws.Range("D4").Value = 0.05740214944097932
$ws.Range("D5").Value = 0.04284130731152905
$ws.Range("D6").Value = 0.0446838664608151
$ws.Range("D7").Value = 0.8545822540618956
$ws.Range("D8").Value = [double]"5.24531075531447e-18"
